$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset all forum-view counts and computed totals/grades to 0 for all
# student rows (rows 2-50, columns B..J), matching the correction
# described in the commit message ("correção das notas do fórum").
$ws.Range("B2:J50").Value = 0
